$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(according to the population census data)" subtitle row.
$ws.Rows.Item(2).Delete()

# Drop the 1989 and 2002 census columns, keeping only the 2014 figures
# (which shift from column D into column B).
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# The remaining rows are now uniformly 20.1pt tall in the refreshed layout.
$ws.Rows("1:5").RowHeight = 20.1
